# Trade #12 closed at 2026-02-16 21:22:35 - leadlag UP +0.000%
# Appends the new trade row (row 11) to the "leadlag" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 11

$ws.Cells.Item($row, 1).Value = 12

# "2026-02-16" looks like a date to Excel's auto-detection, which would
# silently convert it to a date serial + apply a date number format/style.
# Force it to stay literal text, then drop back to the default style so no
# stray style index is left on the cell (matches the rest of the column).
$dateCell = $ws.Cells.Item($row, 2)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-02-16"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 3).Value = "21:22:35"
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "UP"
$ws.Cells.Item($row, 6).Value = 69376.53

# Exit Price is blank (still OPEN) but the column elsewhere stores an
# explicit empty-text cell rather than leaving it fully absent. Assigning a
# bare "" clears/omits the cell entirely, so use a lone apostrophe (Excel's
# "treat as text" quote-prefix) which resolves to an empty string value,
# then reset the style so no quote-prefix style sticks around.
$gCell = $ws.Cells.Item($row, 7)
$gCell.Value = "'"
$gCell.Style = "Normal"

$ws.Cells.Item($row, 8).Value = "OPEN"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0.6506999999999999
$ws.Cells.Item($row, 12).Value = "Coinbase leading with 0.065% move"

# Exit Reason is likewise blank-but-present.
$mCell = $ws.Cells.Item($row, 13)
$mCell.Value = "'"
$mCell.Style = "Normal"

$ws.Cells.Item($row, 14).Value = 0
